# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text updated from "Ready for handoff" to "Handed back: in sync with en-US"
#  - "Latest Target File" (column I) and "Latest Handback File" (column J) are populated
#    for each localized-file row on the zh-cn and de-de sheets, with a hyperlink added on
#    the "Latest Target File" cell pointing back to the original source markdown file.
#  - "Latest Handback DateTime" (column K) is stamped for the de-de sheet (the locale that
#    has actually been handed back in this run).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column updates ---------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

Write-Output "status updated"

# --- Per-file handback details ------------------------------------------------
# Row 2 corresponds to 43bb357f-9e76-4b70-ac86-144c37b4199c.md
# Row 3 corresponds to bcaa9fb7-3425-4373-99c9-edffb332bcce.md
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1d16aea328b0e98670451ce86ee561733788b649/e2e"

$file1Name = "43bb357f-9e76-4b70-ac86-144c37b4199c.md"
$file2Name = "bcaa9fb7-3425-4373-99c9-edffb332bcce.md"
$file1Url  = "$baseUrl/$file1Name"
$file2Url  = "$baseUrl/$file2Name"

$zhXlf1 = "43bb357f-9e76-4b70-ac86-144c37b4199c.97330a2f6bce769c8291f8c4712500ede9b426c1.zh-cn.xlf"
$zhXlf2 = "bcaa9fb7-3425-4373-99c9-edffb332bcce.1450f6c96a2f3cfc0013a8cc94f4e149cec55917.zh-cn.xlf"
$deXlf1 = "43bb357f-9e76-4b70-ac86-144c37b4199c.97330a2f6bce769c8291f8c4712500ede9b426c1.de-de.xlf"
$deXlf2 = "bcaa9fb7-3425-4373-99c9-edffb332bcce.1450f6c96a2f3cfc0013a8cc94f4e149cec55917.de-de.xlf"

$handbackDateTime = "2016-09-03 12:52:05"

# zh-cn sheet: populate Latest Target File (I) with a hyperlink back to the source file,
# and Latest Handback File (J) with the generated xliff file name.
$wsZhCn.Range("I2").Value = $file1Name
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $file1Url, [Type]::Missing, [Type]::Missing, $file1Name) | Out-Null
$wsZhCn.Range("J2").Value = $zhXlf1

$wsZhCn.Range("I3").Value = $file2Name
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $file2Url, [Type]::Missing, [Type]::Missing, $file2Name) | Out-Null
$wsZhCn.Range("J3").Value = $zhXlf2

# de-de sheet: same treatment, plus the Latest Handback DateTime (K) is now known.
$wsDeDe.Range("I2").Value = $file1Name
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $file1Url, [Type]::Missing, [Type]::Missing, $file1Name) | Out-Null
$wsDeDe.Range("J2").Value = $deXlf1
$wsDeDe.Range("K2").Value = $handbackDateTime

$wsDeDe.Range("I3").Value = $file2Name
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $file2Url, [Type]::Missing, [Type]::Missing, $file2Name) | Out-Null
$wsDeDe.Range("J3").Value = $deXlf2
$wsDeDe.Range("K3").Value = $handbackDateTime

Write-Output "handback details written"
